$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
$data = @(
    @(44344, 1, 19, 55.28720246755515),
    @(44345, 7, 23, 66.92661351335623),
    @(44346, 2, 20, 58.19705522900541),
    @(44347, 3, 23, 66.92661351335623),
    @(44348, 1, 23, 66.92661351335623),
    @(44349, 3, 23, 66.92661351335623),
    @(44350, 1, 18, 52.37734970610487),
    @(44351, 4, 21, 61.10690799045568),
    @(44352, 0, 14, 40.73793866030378),
    @(44353, 4, 16, 46.55764418320432),
    @(44354, 5, 18, 52.37734970610487),
    @(44355, 0, 17, 49.4674969446546),
    @(44356, 1, 15, 43.64779142175406),
    @(44357, 3, 17, 49.4674969446546),
    @(44358, 0, 13, 37.82808589885352),
    @(44359, 0, 13, 37.82808589885352),
    @(44360, 0, 9, 26.18867485305244),
    @(44361, 0, 4, 11.63941104580108),
    @(44362, 1, 5, 14.54926380725135),
    @(44363, 0, 4, 11.63941104580108),
    @(44364, 0, 1, 2.90985276145027),
    @(44365, 0, 1, 2.90985276145027),
    @(44366, 1, 2, 5.819705522900541),
    @(44367, 1, 3, 8.729558284350812),
    @(44368, 0, 3, 8.729558284350812),
    @(44369, 0, 2, 5.819705522900541),
    @(44370, 0, 2, 5.819705522900541),
    @(44371, 1, 3, 8.729558284350812),
    @(44372, 0, 3, 8.729558284350812),
    @(44373, 3, 5, 14.54926380725135),
    @(44374, 0, 4, 11.63941104580108),
    @(44375, 1, 5, 14.54926380725135)
)

$startRow = 270
$lastExistingRow = $startRow - 1   # row 269, the last pre-existing data row

# Copy the date-column formatting (style) from the last existing row so the
# newly appended A-column cells carry the same style index (s="2") as the
# rest of the date column.
$ws.Range("A$lastExistingRow").Copy()

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]

    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Range("A$row").PasteSpecial(-4122)   # xlPasteFormats

    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}

$excel.CutCopyMode = $false
